$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New log row (row 6) appended to the GitHub Admin Log sheet.
# Leading apostrophes force text interpretation (so the date-looking string
# and the literal word "False" are stored as text, not a date serial / bool),
# then the style is reset back to Normal so no visible formatting change
# (e.g. quote-prefix) is left on the cell - matching the plain inlineStr
# text cells used by the rest of the sheet.
$ws.Range("A6").Value = "'2025-07-23 12:50:50"
$ws.Range("A6").Style = "Normal"
$ws.Range("B6").Value = "create-team"
$ws.Range("C6").Value = "new-organization97"
$ws.Range("D6").Value = "secondteam"
$ws.Range("E6").Value = "demo"
$ws.Range("F6").Value = "Vignesh2122"
$ws.Range("G6").Value = "pull"
$ws.Range("I6").Value = "'False"
$ws.Range("I6").Style = "Normal"
